$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = 44511
$ws.Cells.Item(2, 10).Value = 500
$ws.Cells.Item(2, 11).Value = 900
$ws.Cells.Item(2, 12).Value = 1000
$ws.Cells.Item(2, 13).Value = 950
$ws.Cells.Item(2, 16).Value = 950

# Row 3
$ws.Cells.Item(3, 4).Value = 44530
$ws.Cells.Item(3, 10).Value = 300

# Row 4
$ws.Cells.Item(4, 4).Value = 44553
$ws.Cells.Item(4, 10).Value = 8000

# Row 5
$ws.Cells.Item(5, 4).Value = 44525
$ws.Cells.Item(5, 10).Value = 360
$ws.Cells.Item(5, 11).Value = 800
$ws.Cells.Item(5, 12).Value = 900
$ws.Cells.Item(5, 13).Value = 850
$ws.Cells.Item(5, 16).Value = 850

# Row 6
$ws.Cells.Item(6, 4).Value = 44537
$ws.Cells.Item(6, 10).Value = 400

# Row 7
$ws.Cells.Item(7, 4).Value = 44505
$ws.Cells.Item(7, 10).Value = 440
$ws.Cells.Item(7, 11).Value = 900
$ws.Cells.Item(7, 12).Value = 1000
$ws.Cells.Item(7, 13).Value = 950
$ws.Cells.Item(7, 16).Value = 950

# Row 8
$ws.Cells.Item(8, 4).Value = 44517
$ws.Cells.Item(8, 10).Value = 500
$ws.Cells.Item(8, 11).Value = 800
$ws.Cells.Item(8, 12).Value = 900
$ws.Cells.Item(8, 13).Value = 850
$ws.Cells.Item(8, 16).Value = 850

# Row 9
$ws.Cells.Item(9, 4).Value = 44518
$ws.Cells.Item(9, 11).Value = 800
$ws.Cells.Item(9, 12).Value = 900
$ws.Cells.Item(9, 13).Value = 850
$ws.Cells.Item(9, 16).Value = 850

# Row 10
$ws.Cells.Item(10, 4).Value = 44503
$ws.Cells.Item(10, 10).Value = 400

# Row 11
$ws.Cells.Item(11, 4).Value = 44532
$ws.Cells.Item(11, 10).Value = 240

# Row 12
$ws.Cells.Item(12, 4).Value = 44545
$ws.Cells.Item(12, 10).Value = 4000

# Row 13
$ws.Cells.Item(13, 4).Value = 44524
$ws.Cells.Item(13, 10).Value = 400

# Row 14
$ws.Cells.Item(14, 4).Value = 44512
$ws.Cells.Item(14, 10).Value = 600

# Row 15
$ws.Cells.Item(15, 4).Value = 44508
$ws.Cells.Item(15, 10).Value = 400
$ws.Cells.Item(15, 11).Value = 900
$ws.Cells.Item(15, 12).Value = 1000
$ws.Cells.Item(15, 13).Value = 950
$ws.Cells.Item(15, 16).Value = 950

# Row 16
$ws.Cells.Item(16, 4).Value = 44523
$ws.Cells.Item(16, 11).Value = 800
$ws.Cells.Item(16, 12).Value = 900
$ws.Cells.Item(16, 13).Value = 850
$ws.Cells.Item(16, 16).Value = 850

# Row 17
$ws.Cells.Item(17, 4).Value = 44476
$ws.Cells.Item(17, 10).Value = 300
$ws.Cells.Item(17, 11).Value = 1100
$ws.Cells.Item(17, 12).Value = 1200
$ws.Cells.Item(17, 13).Value = 1150
$ws.Cells.Item(17, 16).Value = 1150

# Row 18
$ws.Cells.Item(18, 4).Value = 44504
$ws.Cells.Item(18, 10).Value = 500
$ws.Cells.Item(18, 11).Value = 900
$ws.Cells.Item(18, 12).Value = 1000
$ws.Cells.Item(18, 13).Value = 950
$ws.Cells.Item(18, 16).Value = 950

# Row 19
$ws.Cells.Item(19, 4).Value = 44510

# Row 20
$ws.Cells.Item(20, 4).Value = 44516
$ws.Cells.Item(20, 10).Value = 400
$ws.Cells.Item(20, 11).Value = 900
$ws.Cells.Item(20, 12).Value = 1000
$ws.Cells.Item(20, 13).Value = 950
$ws.Cells.Item(20, 16).Value = 950

